# Applies the 4o-mini feedback-sheet text revisions described in the commit
# "most of mini run" to Scleroderma_gen_overall.xlsx.
#
# Sheets: hpi, hist, soc, obj, test (each with header row 1: A=Supports,
# B=Rationale for, C=Against, D=Rationale Against).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# hpi
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("hpi")

$ws.Range("D2").Value = "The lack of skin thickening or tightening is a strong indicator against Scleroderma, as it is a primary feature of the disease."
$ws.Range("D3").Value = "The absence of Raynaud's phenomenon significantly reduces the likelihood of Scleroderma, as it is commonly associated with the condition."
$ws.Range("D4").Value = "A lack of gastrointestinal symptoms suggests that Scleroderma is unlikely, as esophageal involvement is common."
$ws.Range("B5").Value = "Arthralgia (joint pain) and swelling can occur in Scleroderma, indicating systemic involvement."
$ws.Range("D5").Value = "The absence of arthralgia or joint swelling is a strong indicator against Scleroderma, as these symptoms are frequently present."
$ws.Range("B6").Value = "Shortness of breath can indicate interstitial lung disease, a common complication of Scleroderma."
$ws.Range("D6").Value = "The absence of respiratory issues such as shortness of breath or cough suggests that Scleroderma is unlikely, as pulmonary complications are common."

# ---------------------------------------------------------------------------
# hist
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("hist")

$ws.Range("D2").Value = "A lack of any previous skin or connective tissue disorders suggests a lower likelihood of Scleroderma."
$ws.Range("C3").Value = "No history of pulmonary complications"
$ws.Range("D3").Value = "Scleroderma often leads to pulmonary issues; absence of such complications may indicate that Scleroderma is not present."
$ws.Range("B4").Value = "Patients with Scleroderma may have been treated with immunosuppressive medications, indicating a history of autoimmune disease."
$ws.Range("C4").Value = "No previous diagnoses of Raynaud's phenomenon"
$ws.Range("D4").Value = "The absence of Raynaud's phenomenon, which is frequently associated with Scleroderma, suggests a lower likelihood of the disease."
$ws.Range("C5").Value = "No history of joint pain or arthritis"
$ws.Range("D5").Value = "Joint pain or arthritis can be associated with Scleroderma; their absence may indicate that Scleroderma is not present."
$ws.Range("C6").Value = "No use of medications for autoimmune conditions"
$ws.Range("D6").Value = "A lack of treatment for autoimmune conditions suggests that the patient may not have Scleroderma or other related disorders."

# ---------------------------------------------------------------------------
# soc
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("soc")

$ws.Range("B3").Value = "Certain occupations with exposure to silica or other toxins are associated with a higher risk of developing Scleroderma."
$ws.Range("C4").Value = "Non-smoker"
$ws.Range("D4").Value = "Not smoking is associated with a lower risk of developing Scleroderma."
$ws.Range("B5").Value = "Chronic stress is thought to play a role in the development of autoimmune diseases, including Scleroderma."
$ws.Range("D5").Value = "Participation in activities that reduce stress may lower the risk of developing autoimmune diseases, including Scleroderma."
$ws.Range("A6").Value = "History of other autoimmune conditions"
$ws.Range("B6").Value = "Having other autoimmune conditions increases the likelihood of developing Scleroderma due to shared genetic and environmental factors."
$ws.Range("D6").Value = "Absence of other autoimmune diseases suggests a lower likelihood of developing Scleroderma."

# ---------------------------------------------------------------------------
# obj
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("obj")

$ws.Range("D2").Value = "Normal findings in skin texture and elasticity suggest that Scleroderma is unlikely, as skin changes are a key feature."
$ws.Range("C3").Value = "Absence of digital ulcers or skin lesions"
$ws.Range("D3").Value = "The absence of these findings makes Scleroderma less likely, as they are common in affected individuals."
$ws.Range("B4").Value = "The presence of telangiectasia is often seen in Scleroderma and is indicative of vascular changes associated with the disease."
$ws.Range("C4").Value = "Normal capillary refill time"
$ws.Range("D4").Value = "A normal capillary refill time indicates good peripheral circulation, which is inconsistent with Scleroderma."
$ws.Range("B5").Value = "Digital ulcers are frequently observed in Scleroderma patients due to poor circulation and skin changes."
$ws.Range("C5").Value = "No signs of pulmonary hypertension on auscultation"
$ws.Range("D5").Value = "The absence of abnormal lung sounds or signs of pulmonary hypertension suggests that Scleroderma is not present."
$ws.Range("B6").Value = "Pulmonary hypertension can develop in Scleroderma patients, and its detection during a physical exam supports the diagnosis."
$ws.Range("C6").Value = "Normal blood pressure and heart rate"
$ws.Range("D6").Value = "Normal vital signs indicate that there is no significant cardiovascular involvement, which is often seen in Scleroderma."

# ---------------------------------------------------------------------------
# test
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("test")

$ws.Range("D3").Value = "Normal pulmonary function tests suggest that there is no lung involvement, which is often seen in Scleroderma."
$ws.Range("B4").Value = "Restrictive lung disease is a common pulmonary complication of Scleroderma, indicating involvement of the lungs."
$ws.Range("D4").Value = "A negative ANA test is a strong indicator against the presence of autoimmune diseases like Scleroderma."
$ws.Range("B5").Value = "A skin biopsy revealing collagen deposition is indicative of the fibrotic changes seen in Scleroderma."
$ws.Range("D5").Value = "A normal skin biopsy indicates the absence of the fibrotic changes characteristic of Scleroderma."
$ws.Range("B6").Value = "Decreased peristalsis in the esophagus is a common gastrointestinal manifestation of Scleroderma."
$ws.Range("D6").Value = "Normal esophageal motility studies suggest that there is no gastrointestinal involvement, which is common in Scleroderma."
